$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44637
$ws.Range("J2").Value = 38
$ws.Range("K2").Value = 25000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 25000
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Provincia de Quillota"
$ws.Range("P2").Value = 1000

$ws.Range("D3").Value = 44601
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 23000
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = 23600
$ws.Range("N3").Value = "`$/malla 25 kilos"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 944

$ws.Range("D4").Value = 44848
$ws.Range("J4").Value = 38
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 30000
$ws.Range("N4").Value = "`$/malla 25 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 1200

$ws.Range("D5").Value = 44525
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16479
$ws.Range("N5").Value = "`$/malla 25 kilos"
$ws.Range("O5").Value = "Provincia de Talca"
$ws.Range("P5").Value = 659

$ws.Range("D6").Value = 44483
$ws.Range("J6").Value = 55
$ws.Range("K6").Value = 29000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29455
$ws.Range("N6").Value = "`$/malla 25 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 1178

$ws.Range("D7").Value = 44524
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16538
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 662

$ws.Range("D8").Value = 44250
$ws.Range("J8").Value = 38
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 18000
$ws.Range("N8").Value = "`$/malla 25 kilos"
$ws.Range("O8").Value = "Provincia de Talca"
$ws.Range("P8").Value = 720

$ws.Range("D9").Value = 44543
$ws.Range("J9").Value = 35
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 720

$ws.Range("D10").Value = 44399
$ws.Range("J10").Value = 38
$ws.Range("K10").Value = 33000
$ws.Range("L10").Value = 33000
$ws.Range("M10").Value = 33000
$ws.Range("N10").Value = "`$/malla 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1320

$ws.Range("D11").Value = 44677
$ws.Range("J11").Value = 65
$ws.Range("K11").Value = 22000
$ws.Range("L11").Value = 23000
$ws.Range("M11").Value = 22462
$ws.Range("N11").Value = "`$/malla 25 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 898

$ws.Range("D12").Value = 44859
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 24000
$ws.Range("M12").Value = 24000
$ws.Range("N12").Value = "`$/malla 25 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 960

$ws.Range("D13").Value = 44473
$ws.Range("J13").Value = 85
$ws.Range("K13").Value = 35000
$ws.Range("L13").Value = 36000
$ws.Range("M13").Value = 35471
$ws.Range("N13").Value = "`$/malla 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1419

$ws.Range("D14").Value = 44515
$ws.Range("J14").Value = 73
$ws.Range("K14").Value = 16000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16521
$ws.Range("N14").Value = "`$/saco 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 661

$ws.Range("D15").Value = 44253
$ws.Range("J15").Value = 38
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 18000
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Talca"
$ws.Range("P15").Value = 720

$ws.Range("D16").Value = 44372
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 33000
$ws.Range("L16").Value = 34000
$ws.Range("M16").Value = 33500
$ws.Range("N16").Value = "`$/saco 25 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 1340

$ws.Range("D17").Value = 44900
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21479
$ws.Range("N17").Value = "`$/saco 25 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 859

$ws.Range("D18").Value = 44523
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 16500
$ws.Range("M18").Value = 16250
$ws.Range("N18").Value = "`$/malla 25 kilos"
$ws.Range("O18").Value = "Provincia de Talca"
$ws.Range("P18").Value = 650

$ws.Range("D19").Value = 44628
$ws.Range("J19").Value = 73
$ws.Range("K19").Value = 23000
$ws.Range("L19").Value = 24000
$ws.Range("M19").Value = 23521
$ws.Range("N19").Value = "`$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Quillota"
$ws.Range("P19").Value = 941

$ws.Range("D20").Value = 44537
$ws.Range("J20").Value = 78
$ws.Range("K20").Value = 27000
$ws.Range("L20").Value = 28000
$ws.Range("M20").Value = 27487
$ws.Range("N20").Value = "`$/malla 25 kilos"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 1099

$ws.Range("D21").Value = 44874
$ws.Range("J21").Value = 78
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 21000
$ws.Range("M21").Value = 20513
$ws.Range("N21").Value = "`$/malla 25 kilos"
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 821

$ws.Range("D22").Value = 44858
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 24500
$ws.Range("N22").Value = "`$/malla 25 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 980

$ws.Range("D23").Value = 44165
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 22000
$ws.Range("N23").Value = "`$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Quillota"
$ws.Range("P23").Value = 880

$ws.Range("D24").Value = 44365
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 23000
$ws.Range("M24").Value = 22500
$ws.Range("N24").Value = "`$/malla 25 kilos"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 900

$ws.Range("D25").Value = 44510
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = 16500
$ws.Range("L25").Value = 17000
$ws.Range("M25").Value = 16740
$ws.Range("N25").Value = "`$/saco 25 kilos"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 670

$ws.Range("D26").Value = 44550
$ws.Range("J26").Value = 73
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17521
$ws.Range("N26").Value = "`$/saco 25 kilos"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 701

$ws.Range("D27").Value = 44876
$ws.Range("J27").Value = 38
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("N27").Value = "`$/malla 25 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 720

$ws.Range("D28").Value = 44411
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 34000
$ws.Range("L28").Value = 34000
$ws.Range("M28").Value = 34000
$ws.Range("N28").Value = "`$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1360

$ws.Range("D29").Value = 44901
$ws.Range("J29").Value = 65
$ws.Range("K29").Value = 18000
$ws.Range("L29").Value = 19000
$ws.Range("M29").Value = 18462
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 738

$ws.Range("D30").Value = 44159
$ws.Range("J30").Value = 35
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 22000
$ws.Range("M30").Value = 22000
$ws.Range("N30").Value = "`$/malla 25 kilos"
$ws.Range("O30").Value = "Provincia de Quillota"
$ws.Range("P30").Value = 880

$ws.Range("D31").Value = 44546
$ws.Range("J31").Value = 75
$ws.Range("K31").Value = 18000
$ws.Range("L31").Value = 18500
$ws.Range("M31").Value = 18267
$ws.Range("N31").Value = "`$/saco 25 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 731

$ws.Range("D32").Value = 44160
$ws.Range("J32").Value = 35
$ws.Range("K32").Value = 21000
$ws.Range("L32").Value = 21000
$ws.Range("M32").Value = 21000
$ws.Range("N32").Value = "`$/saco 25 kilos"
$ws.Range("O32").Value = "Provincia de Quillota"
$ws.Range("P32").Value = 840

$ws.Range("D33").Value = 44484
$ws.Range("J33").Value = 71
$ws.Range("K33").Value = 29000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = 29507
$ws.Range("N33").Value = "`$/saco 25 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1180

$ws.Range("D34").Value = 44526
$ws.Range("J34").Value = 73
$ws.Range("K34").Value = 16000
$ws.Range("L34").Value = 17000
$ws.Range("M34").Value = 16521
$ws.Range("N34").Value = "`$/saco 25 kilos"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 661

$ws.Range("D35").Value = 44370
$ws.Range("J35").Value = 45
$ws.Range("K35").Value = 32000
$ws.Range("L35").Value = 32000
$ws.Range("M35").Value = 32000
$ws.Range("N35").Value = "`$/malla 25 kilos"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 1280

$ws.Range("D36").Value = 44567
$ws.Range("J36").Value = 68
$ws.Range("K36").Value = 24000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 24559
$ws.Range("N36").Value = "`$/malla 25 kilos"
$ws.Range("O36").Value = "Provincia de Limarí"
$ws.Range("P36").Value = 982

$ws.Range("D37").Value = 44343
$ws.Range("J37").Value = 40
$ws.Range("K37").Value = 28000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 28000
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 1120

$ws.Range("D38").Value = 44676
$ws.Range("J38").Value = 73
$ws.Range("K38").Value = 23000
$ws.Range("L38").Value = 24000
$ws.Range("M38").Value = 23479
$ws.Range("N38").Value = "`$/malla 25 kilos"
$ws.Range("O38").Value = "Provincia de Limarí"
$ws.Range("P38").Value = 939

$ws.Range("D39").Value = 44487
$ws.Range("J39").Value = 73
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 21000
$ws.Range("M39").Value = 20521
$ws.Range("N39").Value = "`$/malla 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 821

$ws.Range("D40").Value = 44875
$ws.Range("J40").Value = 68
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 16000
$ws.Range("M40").Value = 15559
$ws.Range("N40").Value = "`$/malla 25 kilos"
$ws.Range("O40").Value = "Provincia de Quillota"
$ws.Range("P40").Value = 622

$ws.Range("D41").Value = 44868
$ws.Range("J41").Value = 76
$ws.Range("K41").Value = 22000
$ws.Range("L41").Value = 23000
$ws.Range("M41").Value = 22500
$ws.Range("N41").Value = "`$/malla 25 kilos"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 900

$ws.Range("D42").Value = 44161
$ws.Range("J42").Value = 35
$ws.Range("K42").Value = 21000
$ws.Range("L42").Value = 21000
$ws.Range("M42").Value = 21000
$ws.Range("N42").Value = "`$/saco 25 kilos"
$ws.Range("O42").Value = "Provincia de Quillota"
$ws.Range("P42").Value = 840

$ws.Range("D43").Value = 44252
$ws.Range("J43").Value = 40
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 19000
$ws.Range("M43").Value = 18625
$ws.Range("N43").Value = "`$/malla 25 kilos"
$ws.Range("O43").Value = "Provincia de Talca"
$ws.Range("P43").Value = 745

$ws.Range("D44").Value = 44475
$ws.Range("J44").Value = 73
$ws.Range("K44").Value = 25000
$ws.Range("L44").Value = 26000
$ws.Range("M44").Value = 25479
$ws.Range("N44").Value = "`$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 1019

$ws.Range("D45").Value = 44469
$ws.Range("J45").Value = 73
$ws.Range("K45").Value = 28000
$ws.Range("L45").Value = 29000
$ws.Range("M45").Value = 28521
$ws.Range("N45").Value = "`$/malla 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 1141

$ws.Range("D46").Value = 44410
$ws.Range("J46").Value = 35
$ws.Range("K46").Value = 34000
$ws.Range("L46").Value = 34000
$ws.Range("M46").Value = 34000
$ws.Range("N46").Value = "`$/malla 25 kilos"
$ws.Range("O46").Value = "Provincia de Limarí"
$ws.Range("P46").Value = 1360

$ws.Range("D47").Value = 44509
$ws.Range("J47").Value = 80
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 16000
$ws.Range("M47").Value = 15500
$ws.Range("N47").Value = "`$/malla 25 kilos"
$ws.Range("O47").Value = "Provincia de Limarí"
$ws.Range("P47").Value = 620

$ws.Range("D48").Value = 44560
$ws.Range("J48").Value = 50
$ws.Range("K48").Value = 27000
$ws.Range("L48").Value = 28000
$ws.Range("M48").Value = 27500
$ws.Range("N48").Value = "`$/malla 25 kilos"
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 1100

$ws.Range("D49").Value = 44481
$ws.Range("J49").Value = 63
$ws.Range("K49").Value = 22000
$ws.Range("L49").Value = 23000
$ws.Range("M49").Value = 22476
$ws.Range("N49").Value = "`$/saco 25 kilos"
$ws.Range("O49").Value = "Provincia de Limarí"
$ws.Range("P49").Value = 899

$ws.Range("D50").Value = 44476
$ws.Range("J50").Value = 73
$ws.Range("K50").Value = 23000
$ws.Range("L50").Value = 24000
$ws.Range("M50").Value = 23521
$ws.Range("N50").Value = "`$/saco 25 kilos"
$ws.Range("O50").Value = "Provincia de Limarí"
$ws.Range("P50").Value = 941

$ws.Range("D51").Value = 44508
$ws.Range("J51").Value = 68
$ws.Range("K51").Value = 16000
$ws.Range("L51").Value = 17000
$ws.Range("M51").Value = 16515
$ws.Range("N51").Value = "`$/malla 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 661

$ws.Range("D52").Value = 44894
$ws.Range("J52").Value = 73
$ws.Range("K52").Value = 22000
$ws.Range("L52").Value = 23000
$ws.Range("M52").Value = 22521
$ws.Range("N52").Value = "`$/saco 25 kilos"
$ws.Range("O52").Value = "Región Metropolitana"
$ws.Range("P52").Value = 901

$ws.Range("D53").Value = 44448
$ws.Range("J53").Value = 45
$ws.Range("K53").Value = 32000
$ws.Range("L53").Value = 32000
$ws.Range("M53").Value = 32000
$ws.Range("N53").Value = "`$/malla 25 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 1280

$ws.Range("D54").Value = 44529
$ws.Range("J54").Value = 73
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = 17521
$ws.Range("N54").Value = "`$/saco 25 kilos"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 701

$ws.Range("D55").Value = 44908
$ws.Range("J55").Value = 20
$ws.Range("K55").Value = 18000
$ws.Range("L55").Value = 18000
$ws.Range("M55").Value = 18000
$ws.Range("N55").Value = "`$/malla 25 kilos"
$ws.Range("O55").Value = "Provincia de Limarí"
$ws.Range("P55").Value = 720

$ws.Range("D56").Value = 44532
$ws.Range("J56").Value = 73
$ws.Range("K56").Value = 28000
$ws.Range("L56").Value = 29000
$ws.Range("M56").Value = 28521
$ws.Range("N56").Value = "`$/saco 25 kilos"
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 1141

$ws.Range("D57").Value = 44767
$ws.Range("J57").Value = 45
$ws.Range("K57").Value = 37000
$ws.Range("L57").Value = 38000
$ws.Range("M57").Value = 37556
$ws.Range("N57").Value = "`$/saco 25 kilos"
$ws.Range("O57").Value = "Provincia de Limarí"
$ws.Range("P57").Value = 1502

$ws.Range("D58").Value = 44634
$ws.Range("J58").Value = 38
$ws.Range("K58").Value = 25000
$ws.Range("L58").Value = 25000
$ws.Range("M58").Value = 25000
$ws.Range("N58").Value = "`$/malla 25 kilos"
$ws.Range("O58").Value = "Provincia de Talca"
$ws.Range("P58").Value = 1000

$ws.Range("D59").Value = 44511
$ws.Range("J59").Value = 73
$ws.Range("K59").Value = 16000
$ws.Range("L59").Value = 17000
$ws.Range("M59").Value = 16479
$ws.Range("N59").Value = "`$/saco 25 kilos"
$ws.Range("O59").Value = "Provincia de Limarí"
$ws.Range("P59").Value = 659

$ws.Range("D60").Value = 44831
$ws.Range("J60").Value = 45
$ws.Range("K60").Value = 28000
$ws.Range("L60").Value = 28000
$ws.Range("M60").Value = 28000
$ws.Range("N60").Value = "`$/saco 25 kilos"
$ws.Range("O60").Value = "Provincia de Limarí"
$ws.Range("P60").Value = 1120

$ws.Range("D61").Value = 44181
$ws.Range("J61").Value = 38
$ws.Range("K61").Value = 26000
$ws.Range("L61").Value = 26000
$ws.Range("M61").Value = 26000
$ws.Range("N61").Value = "`$/malla 25 kilos"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 1040

$ws.Range("D62").Value = 44406
$ws.Range("J62").Value = 35
$ws.Range("K62").Value = 32000
$ws.Range("L62").Value = 32000
$ws.Range("M62").Value = 32000
$ws.Range("N62").Value = "`$/malla 25 kilos"
$ws.Range("O62").Value = "Provincia de Limarí"
$ws.Range("P62").Value = 1280

$ws.Range("D63").Value = 44907
$ws.Range("J63").Value = 73
$ws.Range("K63").Value = 18000
$ws.Range("L63").Value = 19000
$ws.Range("M63").Value = 18521
$ws.Range("N63").Value = "`$/malla 25 kilos"
$ws.Range("O63").Value = "Provincia de Limarí"
$ws.Range("P63").Value = 741

$ws.Range("D64").Value = 44882
$ws.Range("J64").Value = 65
$ws.Range("K64").Value = 19000
$ws.Range("L64").Value = 20000
$ws.Range("M64").Value = 19462
$ws.Range("N64").Value = "`$/saco 25 kilos"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 778

$ws.Range("D65").Value = 44578
$ws.Range("J65").Value = 73
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = 18521
$ws.Range("N65").Value = "`$/malla 25 kilos"
$ws.Range("O65").Value = "Provincia de Talca"
$ws.Range("P65").Value = 741

$ws.Range("D66").Value = 44895
$ws.Range("J66").Value = 73
$ws.Range("K66").Value = 22000
$ws.Range("L66").Value = 23000
$ws.Range("M66").Value = 22521
$ws.Range("N66").Value = "`$/saco 25 kilos"
$ws.Range("O66").Value = "Región Metropolitana"
$ws.Range("P66").Value = 901

$ws.Range("D67").Value = 44452
$ws.Range("J67").Value = 70
$ws.Range("K67").Value = 31000
$ws.Range("L67").Value = 32000
$ws.Range("M67").Value = 31500
$ws.Range("N67").Value = "`$/malla 25 kilos"
$ws.Range("O67").Value = "Provincia de Limarí"
$ws.Range("P67").Value = 1260

$ws.Range("D68").Value = 44376
$ws.Range("J68").Value = 38
$ws.Range("K68").Value = 27000
$ws.Range("L68").Value = 27000
$ws.Range("M68").Value = 27000
$ws.Range("N68").Value = "`$/saco 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 1080

$ws.Range("D69").Value = 44847
$ws.Range("J69").Value = 71
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 31000
$ws.Range("M69").Value = 30493
$ws.Range("N69").Value = "`$/malla 25 kilos"
$ws.Range("O69").Value = "Provincia de Limarí"
$ws.Range("P69").Value = 1220

$ws.Range("D70").Value = 44629
$ws.Range("J70").Value = 45
$ws.Range("K70").Value = 24000
$ws.Range("L70").Value = 25000
$ws.Range("M70").Value = 24444
$ws.Range("N70").Value = "`$/saco 25 kilos"
$ws.Range("O70").Value = "Región Metropolitana"
$ws.Range("P70").Value = 978

$ws.Range("D71").Value = 44536
$ws.Range("J71").Value = 81
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 28000
$ws.Range("M71").Value = 27556
$ws.Range("N71").Value = "`$/saco 25 kilos"
$ws.Range("O71").Value = "Provincia de Limarí"
$ws.Range("P71").Value = 1102

$ws.Range("D72").Value = 44592
$ws.Range("J72").Value = 38
$ws.Range("K72").Value = 22000
$ws.Range("L72").Value = 22000
$ws.Range("M72").Value = 22000
$ws.Range("N72").Value = "`$/malla 25 kilos"
$ws.Range("O72").Value = "Provincia de Talca"
$ws.Range("P72").Value = 880

$ws.Range("D73").Value = 44453
$ws.Range("J73").Value = 73
$ws.Range("K73").Value = 21000
$ws.Range("L73").Value = 22000
$ws.Range("M73").Value = 21521
$ws.Range("N73").Value = "`$/saco 25 kilos"
$ws.Range("O73").Value = "Provincia de Limarí"
$ws.Range("P73").Value = 861

$ws.Range("D74").Value = 44162
$ws.Range("J74").Value = 35
$ws.Range("K74").Value = 17000
$ws.Range("L74").Value = 17000
$ws.Range("M74").Value = 17000
$ws.Range("N74").Value = "`$/saco 25 kilos"
$ws.Range("O74").Value = "Provincia de Quillota"
$ws.Range("P74").Value = 680

$ws.Range("D75").Value = 44468
$ws.Range("J75").Value = 65
$ws.Range("K75").Value = 24000
$ws.Range("L75").Value = 25000
$ws.Range("M75").Value = 24538
$ws.Range("N75").Value = "`$/malla 25 kilos"
$ws.Range("O75").Value = "Provincia de Limarí"
$ws.Range("P75").Value = 982

$ws.Range("D76").Value = 44679
$ws.Range("J76").Value = 77
$ws.Range("K76").Value = 26000
$ws.Range("L76").Value = 27000
$ws.Range("M76").Value = 26506
$ws.Range("N76").Value = "`$/malla 25 kilos"
$ws.Range("O76").Value = "Provincia de Limarí"
$ws.Range("P76").Value = 1060

$ws.Range("D77").Value = 44571
$ws.Range("J77").Value = 73
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = 15479
$ws.Range("N77").Value = "`$/saco 25 kilos"
$ws.Range("O77").Value = "Provincia de Limarí"
$ws.Range("P77").Value = 619

# New row 78
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = "Femacal de La Calera"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44412
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = 5
$ws.Range("F78").Value = 100112022
$ws.Range("G78").Value = "Arveja Verde"
$ws.Range("H78").Value = "Perfection"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 35
$ws.Range("K78").Value = 24000
$ws.Range("L78").Value = 24000
$ws.Range("M78").Value = 24000
$ws.Range("N78").Value = "`$/malla 25 kilos"
$ws.Range("O78").Value = "Provincia de Limarí"
$ws.Range("P78").Value = 960
$ws.Range("Q78").Value = 25
$ws.Range("R78").Value = "Hortaliza"

Write-Output "Done"